# Refresh cryptos list (prices / 1h volume %) -- GitHub Actions scheduled update
# Wed Mar 27 17:42:02 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '69.118.57'
$ws.Range('E2').Value = '  -1.46%  '
$ws.Range('D3').Value = '3.521.79'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.70%  '
$ws.Range('E7').Value = '  -2.44%  '
$ws.Range('D8').Value = '3.516.60'
$ws.Range('E8').Value = '  -1.77%  '
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.184'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.26%  '
$ws.Range('E11').Value = '  -3.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.94'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.70%  '
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.47'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('D15').Value = '4.089.96'
$ws.Range('E15').Value = '  -1.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.01%  '
$ws.Range('D17').Value = '3.513.69'
$ws.Range('E17').Value = '  -1.83%  '
$ws.Range('D18').Value = '69.124.39'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '539.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +13.63%  '
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '20.73'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.11%  '
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.43'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '94.76'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('E28').Value = '  -4.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.18'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.59'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.76%  '
$ws.Range('E31').Value = '  -5.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.70'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.96%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.114'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.48%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '64.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '572.59'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.08'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.11%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.06'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.47%  '
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('D40').Value = '0.0₃0765'
$ws.Range('E40').Value = '  -4.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.134'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.09'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.64%  '
$ws.Range('E43').Value = '  -4.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.54'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.73%  '
$ws.Range('D45').Value = '3.211.26'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0442'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('E48').Value = '  -3.87%  '
$ws.Range('E49').Value = '  -2.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.998'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '136.49'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.72%  '
